$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: record the "Resolved" date for the create_word_xml / unix-mac bug
$ws.Range("D29").Value = 43425
$ws.Range("D29").NumberFormat = "m/d/yy"

# New bug entry in row 32 (row 31 left blank, matching the sheet's existing gaps)
$ws.Range("C32").Value = "code_tree is wrong when created in the vignette with the template/source thing."

# Match the author's last active selection
$ws.Range("C32").Select()
